$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" property value (row 8, column B).
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Update the "Contact" property value (row 10, column B).
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row for the "Jurisdiction" property right after "Contact"
# (before "Description"), pushing everything below down by one row.
$ws.Rows.Item(11).Insert()

# Reuse the formatting of the row above (Contact) for the new row so the
# new cells pick up the existing shared cell style instead of creating a
# brand new one.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
